$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the left (new column A), shifting existing
# columns A:D to B:E. This also shifts the "발행시간/회차/추가주식수/가액"
# header row and all data rows to the right by one column.
$ws.Columns.Item(1).Insert()

# Copy the style of the header cell that landed in B1 (the original
# bold/bordered header style) onto the new A1 header cell so it matches
# the rest of row 1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Fill in the new header and company-name column.
$ws.Range("A1").Value = "회사명"

$company = "에스티팜"
$ws.Range("A2").Value = $company
$ws.Range("A3").Value = $company
$ws.Range("A4").Value = $company
$ws.Range("A5").Value = $company
$ws.Range("A6").Value = $company
$ws.Range("A7").Value = $company
$ws.Range("A8").Value = $company
$ws.Range("A9").Value = $company

# The "회차" column (now column C after the insert) becomes numeric
# instead of text.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 2
